$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes ---
# Stored OOXML width target is 13.4101845877511. This runtime persists
# width as (ColumnWidth + 5/6), so back-solve the ColumnWidth to assign.
$targetColumnWidth = 13.4101845877511 - (5/6)

# Overview sheet: columns E (5) and F (6)
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# zh-cn sheet: column C (3)
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

# de-de sheet: column C (3)
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
